$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing data rows (old rows 31 and 32); this shifts
# the signature block (old rows 37/38) up to rows 35/36 automatically.
$ws.Rows(32).Delete()
$ws.Rows(31).Delete()

# --- Summary header values ---
$ws.Range("E11").Value = 789117
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 8

# --- Data table rows 17-30 (row 16 is unchanged) ---
$data = @(
    @(17, "80054295",   "JORGE LUIS ELJADUE MARTINEZ",     "1707", 29509,  877803),
    @(18, "18389622",   "HERIBERTO AGUIRRE RENDON",        "1707", 29509,  781242),
    @(19, "18386112",   "ELMER PACHON PEREZ",              "1707", 29509,  781242),
    @(20, "80054295",   "JORGE LUIS ELJADUE MARTINEZ",     "1708", 29509,  877803),
    @(21, "18389622",   "HERIBERTO AGUIRRE RENDON",        "1708", 29509,  781242),
    @(22, "18386112",   "ELMER PACHON PEREZ",              "1708", 29509,  781242),
    @(23, "80054295",   "JORGE LUIS ELJADUE MARTINEZ",     "1709", 29509,  877803),
    @(24, "18389622",   "HERIBERTO AGUIRRE RENDON",        "1709", 29509,  781242),
    @(25, "18386112",   "ELMER PACHON PEREZ",              "1709", 29509,  781242),
    @(26, "80054295",   "JORGE LUIS ELJADUE MARTINEZ",     "1710", 29509,  877803),
    @(27, "80054295",   "JORGE LUIS ELJADUE MARTINEZ",     "1711", 29509,  877803),
    @(28, "1193522525", "JEFRIN ANDRES ESTREMOR BLANCO",   "2001", 26400,  900000),
    @(29, "73146547",   "JULIO CESAR PATERNINA FERNANDEZ", "2002", 385084, 25179560),
    @(30, "1050718330", "ALVARO ENRIQUE GULLOSO LEIVA",    "2002", 45680,  1142000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
}
